$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.5
$ws.Range("C5").Value = 14

$ws.Range("B2").Select()
